# feat: add 2022-Q3 data
#
# The workbook originally has two sheets: "总计" (summary) and "2022-Q2"
# (fund holdings detail for that quarter). This edit inserts a new
# "2022-Q3" quarter: the "总计" summary table gets a new top row for
# 2022-Q3 (existing 2022-Q2 row shifts down), a new "2022-Q3" detail
# sheet is created (reusing the original "2022-Q2" sheet's identity so it
# keeps sheetId 2 / stays the 2nd tab), and the original "2022-Q2" detail
# data is preserved by duplicating it into a new trailing "2022-Q2" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Preserve the existing "2022-Q2" detail sheet by copying it to a new
#    sheet placed right after it. The copy keeps all of the old data,
#    formatting and page setup untouched, and becomes the new trailing
#    "2022-Q2" tab.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($null, $oldQ2)
$q2Copy = $wb.Worksheets.Item($oldQ2.Index + 1)

# ---------------------------------------------------------------------
# 2. Turn the ORIGINAL "2022-Q2" sheet into the new "2022-Q3" sheet: wipe
#    its old content and refill it with the Q3 fund holdings, so it keeps
#    its original sheetId / tab position (2nd tab).
# ---------------------------------------------------------------------
$oldQ2.Cells.Clear()
$oldQ2.Name = "2022-Q3"

# Match page margins used by a freshly-authored sheet in this workbook
# (same as the "总计" sheet).
$oldQ2.PageSetup.LeftMargin = 0.75 * 72
$oldQ2.PageSetup.RightMargin = 0.75 * 72
$oldQ2.PageSetup.TopMargin = 1 * 72
$oldQ2.PageSetup.BottomMargin = 1 * 72
$oldQ2.PageSetup.HeaderMargin = 0.5 * 72
$oldQ2.PageSetup.FooterMargin = 0.5 * 72

# Header row - reuse the bold/centered header style from "总计".
$summary = $wb.Worksheets.Item("总计")
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $oldQ2.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$summary.Range("B1:D1").Copy()
$oldQ2.Range("B1:H1").PasteSpecial(-4122)

# Data rows. Columns D-G are stored as text (not numbers) to match the
# existing convention used throughout this workbook.
$rows = @(
    @(0, "001743", "诺安优选回报灵活配置混合", "13.65", "73.02", "4.15", "0.5665", 3),
    @(1, "002319", "大成一带一路灵活配置混合", "1.26", "89.65", "5.76", "0.0726", 5),
    @(2, "003799", "华安新泰利灵活配置混合A", "3.14", "23.43", "0.83", "0.0261", 10),
    @(3, "003800", "华安新泰利灵活配置混合C", "1.65", "23.43", "0.83", "0.0137", 10),
    @(4, "001744", "诺安进取回报灵活配置混合", "0.23", "82.31", "4.00", "0.0092", 8)
)

$textRange = $oldQ2.Range("B2:B6")
$textRange.NumberFormat = "@"
$textRange = $oldQ2.Range("D2:G6")
$textRange.NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $oldQ2.Cells.Item($excelRow, 1).Value = $row[0]
    $oldQ2.Cells.Item($excelRow, 2).Value = $row[1]
    $oldQ2.Cells.Item($excelRow, 3).Value = $row[2]
    $oldQ2.Cells.Item($excelRow, 4).Value = $row[3]
    $oldQ2.Cells.Item($excelRow, 5).Value = $row[4]
    $oldQ2.Cells.Item($excelRow, 6).Value = $row[5]
    $oldQ2.Cells.Item($excelRow, 7).Value = $row[6]
    $oldQ2.Cells.Item($excelRow, 8).Value = $row[7]
}

# A-column cells carry the same style as the "总计" sheet's A2 cell.
$summary.Range("A2").Copy()
$oldQ2.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Rename the duplicated sheet back to "2022-Q2" now that the original
#    name has been freed up.
# ---------------------------------------------------------------------
$q2Copy.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 4. Update the "总计" summary sheet: insert a new row for 2022-Q3 above
#    the existing 2022-Q2 row.
# ---------------------------------------------------------------------
$summary.Range("A3").EntireRow.Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.69

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 4
$summary.Cells.Item(3, 4).Value = 0.36

# The row-insert above leaves the shifted A3 cell with a slightly
# different (borderless) auto-derived style - restore the original A2
# style (used uniformly by column A in this sheet).
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
